$d = $word.ActiveDocument

# --- "Programa resumido" paragraph ---
# Originally a single run of text with the 7 numbered topics concatenated.
# Split it into separate runs/text-parts joined by manual line breaks (<w:br/>)
# right after each topic's trailing colon / period, i.e. before the next
# item's leading digit.
$d.Content.Find.Execute("Introdução:2-", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Introdução:^l2-", 2) | Out-Null
$d.Content.Find.Execute("difusão:3-", $false, $false, $false, $false, $false, `
    $true, 1, $false, "difusão:^l3-", 2) | Out-Null
$d.Content.Find.Execute("fluxos:4 -", $false, $false, $false, $false, $false, `
    $true, 1, $false, "fluxos:^l4 -", 2) | Out-Null
$d.Content.Find.Execute("de massa:5-", $false, $false, $false, $false, $false, `
    $true, 1, $false, "de massa:^l5-", 2) | Out-Null
$d.Content.Find.Execute("química:6-", $false, $false, $false, $false, $false, `
    $true, 1, $false, "química:^l6-", 2) | Out-Null
$d.Content.Find.Execute("química:7-", $false, $false, $false, $false, $false, `
    $true, 1, $false, "química:^l7-", 2) | Out-Null

# --- "Critério" paragraph ---
# Insert a manual line break between the formula "...)/3" and the following
# legend "P2 = Nota da Prova...".
$d.Content.Find.Execute("2*P2)/3P2 =", $false, $false, $false, $false, $false, `
    $true, 1, $false, "2*P2)/3^lP2 =", 2) | Out-Null

Write-Output "Edit applied"
